# error solve ifrs list
# Update D2:AJ6 with corrected values, and clear out the stray/duplicate
# data that had leaked into rows 7-9 (columns D onward), keeping only
# the row number (A), company group (B) and company name (C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# --- Row 2 ---
$ws.Range("D2").Value = 2451
$ws.Range("E2").Value = 184
$ws.Range("F2").Value = 184
$ws.Range("G2").Value = 218
$ws.Range("H2").Value = 181
$ws.Range("I2").Value = 186
$ws.Range("J2").Value = -5
$ws.Range("K2").Value = 1840
$ws.Range("L2").Value = 460
$ws.Range("M2").Value = 1380
$ws.Range("N2").Value = 1315
$ws.Range("O2").Value = 65
$ws.Range("P2").Value = 100
$ws.Range("Q2").Value = 85
$ws.Range("R2").Value = -71
$ws.Range("S2").Value = 27
$ws.Range("T2").Value = 13
$ws.Range("U2").Value = 71
$ws.Range("V2").Value = 100
$ws.Range("W2").Value = 7.49
$ws.Range("X2").Value = 7.38
$ws.Range("Y2").Value = 15.05
$ws.Range("Z2").Value = 10.67
$ws.Range("AA2").Value = 33.32
$ws.Range("AB2").Value = 1216.01
$ws.Range("AC2").Value = 6492
$ws.Range("AD2").Value = 4.68
$ws.Range("AE2").Value = 46038
$ws.Range("AF2").Value = 0.66
$ws.Range("AG2").Value = 1400
$ws.Range("AH2").Value = 4.61
$ws.Range("AI2").Value = 21.56
$ws.Range("AJ2").Value = 2857223

# --- Row 3 ---
$ws.Range("D3").Value = 2313
$ws.Range("E3").Value = 189
$ws.Range("F3").Value = 189
$ws.Range("G3").Value = 221
$ws.Range("H3").Value = 158
$ws.Range("I3").Value = 157
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1953
$ws.Range("L3").Value = 437
$ws.Range("M3").Value = 1516
$ws.Range("N3").Value = 1437
$ws.Range("O3").Value = 79
$ws.Range("P3").Value = 100
$ws.Range("Q3").Value = 180
$ws.Range("R3").Value = -146
$ws.Range("S3").Value = -31
$ws.Range("T3").Value = 25
$ws.Range("U3").Value = 154
$ws.Range("V3").Value = 116
$ws.Range("W3").Value = 8.19
$ws.Range("X3").Value = 6.81
$ws.Range("Y3").Value = 11.42
$ws.Range("Z3").Value = 8.31
$ws.Range("AA3").Value = 28.85
$ws.Range("AB3").Value = 1330.36
$ws.Range("AC3").Value = 5501
$ws.Range("AD3").Value = 5.45
$ws.Range("AE3").Value = 50292
$ws.Range("AF3").Value = 0.6
$ws.Range("AG3").Value = 910
$ws.Range("AH3").Value = 3.03
$ws.Range("AI3").Value = 16.54
$ws.Range("AJ3").Value = 2857223

# --- Row 4 ---
$ws.Range("D4").Value = 2134
$ws.Range("E4").Value = 165
$ws.Range("F4").Value = 165
$ws.Range("G4").Value = 186
$ws.Range("H4").Value = 140
$ws.Range("I4").Value = 139
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2024
$ws.Range("L4").Value = 392
$ws.Range("M4").Value = 1632
$ws.Range("N4").Value = 1553
$ws.Range("O4").Value = 79
$ws.Range("P4").Value = 100
$ws.Range("Q4").Value = 24
$ws.Range("R4").Value = 95
$ws.Range("S4").Value = -109
$ws.Range("T4").Value = 16
$ws.Range("U4").Value = 9
$ws.Range("V4").Value = 34
$ws.Range("W4").Value = 7.72
$ws.Range("X4").Value = 6.55
$ws.Range("Y4").Value = 9.33
$ws.Range("Z4").Value = 7.03
$ws.Range("AA4").Value = 24.04
$ws.Range("AB4").Value = 1444.29
$ws.Range("AC4").Value = 4880
$ws.Range("AD4").Value = 6.97
$ws.Range("AE4").Value = 54362
$ws.Range("AF4").Value = 0.63
$ws.Range("AG4").Value = 910
$ws.Range("AH4").Value = 2.68
$ws.Range("AI4").Value = 18.65
$ws.Range("AJ4").Value = 2857223

# --- Row 5 ---
$ws.Range("D5").Value = 2576
$ws.Range("E5").Value = 189
$ws.Range("F5").Value = 189
$ws.Range("G5").Value = 181
$ws.Range("H5").Value = 127
$ws.Range("I5").Value = 123
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 2170
$ws.Range("L5").Value = 453
$ws.Range("M5").Value = 1717
$ws.Range("N5").Value = 1639
$ws.Range("O5").Value = 78
$ws.Range("P5").Value = 100
$ws.Range("Q5").Value = -3
$ws.Range("R5").Value = 61
$ws.Range("S5").Value = -1
$ws.Range("T5").Value = 32
$ws.Range("U5").Value = -35
$ws.Range("V5").Value = 57
$ws.Range("W5").Value = 7.32
$ws.Range("X5").Value = 4.93
$ws.Range("Y5").Value = 7.69
$ws.Range("Z5").Value = 6.06
$ws.Range("AA5").Value = 26.38
$ws.Range("AB5").Value = 1543.61
$ws.Range("AC5").Value = 4295
$ws.Range("AD5").Value = 9.6
$ws.Range("AE5").Value = 57362
$ws.Range("AF5").Value = 0.72
$ws.Range("AG5").Value = 700
$ws.Range("AH5").Value = 1.7
$ws.Range("AI5").Value = 16.3
$ws.Range("AJ5").Value = 2857223

# --- Row 6 (note: J6 and O6 are intentionally absent, as in the source) ---
$ws.Range("D6").Value = 3205
$ws.Range("E6").Value = 139
$ws.Range("F6").Value = 139
$ws.Range("G6").Value = 167
$ws.Range("H6").Value = 125
$ws.Range("I6").Value = 120
$ws.Range("K6").Value = 2290
$ws.Range("L6").Value = 477
$ws.Range("M6").Value = 1813
$ws.Range("N6").Value = 1740
$ws.Range("P6").Value = 100
$ws.Range("Q6").Value = -33
$ws.Range("R6").Value = 51
$ws.Range("S6").Value = -27
$ws.Range("T6").Value = 47
$ws.Range("U6").Value = -80
$ws.Range("V6").Value = 60
$ws.Range("W6").Value = 4.34
$ws.Range("X6").Value = 3.9
$ws.Range("Y6").Value = 7.11
$ws.Range("Z6").Value = 5.6
$ws.Range("AA6").Value = 26.28
$ws.Range("AB6").Value = 1644.72
$ws.Range("AC6").Value = 4201
$ws.Range("AD6").Value = 7.12
$ws.Range("AE6").Value = 61107
$ws.Range("AF6").Value = 0.49
$ws.Range("AG6").Value = 700
$ws.Range("AH6").Value = 2.34
$ws.Range("AI6").Value = 16.6
$ws.Range("AJ6").Value = 2857223

# --- Rows 7-9: clear the stray data columns (D onward), keep A/B/C ---
$ws.Range("D7:AJ9").ClearContents()
